$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E44").Value = 70
$ws.Range("H44").Value = 20
$ws.Range("K44").Value = 0
$ws.Range("N44").Value = 50
$ws.Range("Q44").Value = 40
$ws.Range("T44").Value = 80
$ws.Range("W44").Value = 60
$ws.Range("Z44").Value = 100
$ws.Range("AC44").Value = 30

$excel.CalculateFull()
